$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 151; existing rows 151-198 shift down to 152-199.
$ws.Rows(151).Insert()

# Populate the newly inserted row 151 with the new weekly price-report entry.
$ws.Cells.Item(151, 1).Value = 5
$ws.Cells.Item(151, 2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item(151, 3).Value = 'Maule'
$ws.Cells.Item(151, 4).Value = 44988
$ws.Cells.Item(151, 5).Value = 7
$ws.Cells.Item(151, 6).Value = 100112030
$ws.Cells.Item(151, 7).Value = 'Poroto granado'
$ws.Cells.Item(151, 8).Value = 'Sin especificar'
$ws.Cells.Item(151, 9).Value = 'Primera'
$ws.Cells.Item(151, 10).Value = 400
$ws.Cells.Item(151, 11).Value = 30000
$ws.Cells.Item(151, 12).Value = 30000
$ws.Cells.Item(151, 13).Value = 30000
$ws.Cells.Item(151, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(151, 15).Value = 'Región del Maule'
$ws.Cells.Item(151, 16).Value = 1200
$ws.Cells.Item(151, 17).Value = 25
$ws.Cells.Item(151, 18).Value = 'Hortaliza'
